# "Update to latest indices and docs"
#
# The underlying model was re-run; refresh the raw output sheets
# ("raw East_perf" / "raw West_perf") with the latest performance
# indices. The "For publication" sheet reads these via live formulas
# (e.g. ='raw East_perf'!B1), so once the raw values change and the
# workbook recalculates, the published summary table picks up the new
# numbers automatically.

$wb = $excel.ActiveWorkbook

$east = $wb.Worksheets.Item("raw East_perf")
$west = $wb.Worksheets.Item("raw West_perf")

$eastArr = New-Object 'object[,]' 6,13
$eastArr[0,0] = 0
$eastArr[0,1] = 0
$eastArr[0,2] = 0
$eastArr[0,3] = 0.45
$eastArr[0,4] = 0.67400000000000004
$eastArr[0,5] = 0.82099999999999995
$eastArr[0,6] = 0.32200000000000001
$eastArr[0,7] = 1
$eastArr[0,8] = 1
$eastArr[0,9] = 0
$eastArr[0,10] = 0
$eastArr[0,11] = 100
$eastArr[0,12] = 1.7030000000000001
$eastArr[1,0] = 5.8970000000000002
$eastArr[1,1] = 5.8970000000000002
$eastArr[1,2] = 5.8970000000000002
$eastArr[1,3] = 0.42
$eastArr[1,4] = 0.59499999999999997
$eastArr[1,5] = 0.71
$eastArr[1,6] = 0.316
$eastArr[1,7] = 0.85
$eastArr[1,8] = 0.84399999999999997
$eastArr[1,9] = 0
$eastArr[1,10] = 0
$eastArr[1,11] = 100
$eastArr[1,12] = 0
$eastArr[2,0] = 11.792999999999999
$eastArr[2,1] = 11.792999999999999
$eastArr[2,2] = 11.792999999999999
$eastArr[2,3] = 0.39
$eastArr[2,4] = 0.51500000000000001
$eastArr[2,5] = 0.59299999999999997
$eastArr[2,6] = 0.31
$eastArr[2,7] = 0.68899999999999995
$eastArr[2,8] = 0.67900000000000005
$eastArr[2,9] = 0
$eastArr[2,10] = 0
$eastArr[2,11] = 100
$eastArr[2,12] = 0
$eastArr[3,0] = 17.690000000000001
$eastArr[3,1] = 17.690000000000001
$eastArr[3,2] = 17.690000000000001
$eastArr[3,3] = 0.36
$eastArr[3,4] = 0.434
$eastArr[3,5] = 0.47099999999999997
$eastArr[3,6] = 0.3
$eastArr[3,7] = 0.50800000000000001
$eastArr[3,8] = 0.495
$eastArr[3,9] = 0
$eastArr[3,10] = 0.002
$eastArr[3,11] = 99.91
$eastArr[3,12] = 0
$eastArr[4,0] = 12.279
$eastArr[4,1] = 15.244999999999999
$eastArr[4,2] = 20.606999999999999
$eastArr[4,3] = 0.38900000000000001
$eastArr[4,4] = 0.501
$eastArr[4,5] = 0.53
$eastArr[4,6] = 0.23300000000000001
$eastArr[4,7] = 0.27600000000000002
$eastArr[4,8] = 0.27600000000000002
$eastArr[4,9] = 0
$eastArr[4,10] = 0
$eastArr[4,11] = 100
$eastArr[4,12] = 2.238
$eastArr[5,0] = 12.087
$eastArr[5,1] = 13.739000000000001
$eastArr[5,2] = 17.321999999999999
$eastArr[5,3] = 0.38900000000000001
$eastArr[5,4] = 0.50700000000000001
$eastArr[5,5] = 0.55700000000000005
$eastArr[5,6] = 0.29799999999999999
$eastArr[5,7] = 0.438
$eastArr[5,8] = 0.438
$eastArr[5,9] = 0
$eastArr[5,10] = 0
$eastArr[5,11] = 100
$eastArr[5,12] = 1.6519999999999999

$westArr = New-Object 'object[,]' 6,13
$westArr[0,0] = 0
$westArr[0,1] = 0
$westArr[0,2] = 0
$westArr[0,3] = 0.51500000000000001
$westArr[0,4] = 1.093
$westArr[0,5] = 1.081
$westArr[0,6] = 0.377
$westArr[0,7] = 1
$westArr[0,8] = 1
$westArr[0,9] = 0
$westArr[0,10] = 0
$westArr[0,11] = 100
$westArr[0,12] = 1.7030000000000001
$westArr[1,0] = 1.8120000000000001
$westArr[1,1] = 1.8140000000000001
$westArr[1,2] = 1.8140000000000001
$westArr[1,3] = 0.46300000000000002
$westArr[1,4] = 0.92700000000000005
$westArr[1,5] = 0.876
$westArr[1,6] = 0.34899999999999998
$westArr[1,7] = 0.77900000000000003
$westArr[1,8] = 0.77100000000000002
$westArr[1,9] = 0
$westArr[1,10] = 0
$westArr[1,11] = 100
$westArr[1,12] = 0.012999999999999999
$westArr[2,0] = 3.419
$westArr[2,1] = 3.5979999999999999
$westArr[2,2] = 3.5840000000000001
$westArr[2,3] = 0.41899999999999998
$westArr[2,4] = 0.77200000000000002
$westArr[2,5] = 0.68100000000000005
$westArr[2,6] = 0.318
$westArr[2,7] = 0.56000000000000005
$westArr[2,8] = 0.54800000000000004
$westArr[2,9] = 0
$westArr[2,10] = 0
$westArr[2,11] = 100
$westArr[2,12] = 0.27
$westArr[3,0] = 4.8310000000000004
$westArr[3,1] = 5.0229999999999997
$westArr[3,2] = 4.9489999999999998
$westArr[3,3] = 0.378
$westArr[3,4] = 0.627
$westArr[3,5] = 0.501
$westArr[3,6] = 0.25700000000000001
$westArr[3,7] = 0.34399999999999997
$westArr[3,8] = 0.33300000000000002
$westArr[3,9] = 0
$westArr[3,10] = 0.0040000000000000001
$westArr[3,11] = 99.82
$westArr[3,12] = 0.33500000000000002
$westArr[4,0] = 3.3719999999999999
$westArr[4,1] = 3.7919999999999998
$westArr[4,2] = 4.4619999999999997
$westArr[4,3] = 0.41899999999999998
$westArr[4,4] = 0.76200000000000001
$westArr[4,5] = 0.628
$westArr[4,6] = 0.215
$westArr[4,7] = 0.24
$westArr[4,8] = 0.23899999999999999
$westArr[4,9] = 0
$westArr[4,10] = 0
$westArr[4,11] = 100
$westArr[4,12] = 1.6990000000000001
$westArr[5,0] = 3.25
$westArr[5,1] = 3.4340000000000002
$westArr[5,2] = 3.851
$westArr[5,3] = 0.42099999999999999
$westArr[5,4] = 0.77600000000000002
$westArr[5,5] = 0.66600000000000004
$westArr[5,6] = 0.28599999999999998
$westArr[5,7] = 0.38400000000000001
$westArr[5,8] = 0.38100000000000001
$westArr[5,9] = 0
$westArr[5,10] = 0
$westArr[5,11] = 100
$westArr[5,12] = 1.177

$east.Range("B2:N7").Value2 = $eastArr
$west.Range("B2:N7").Value2 = $westArr

$excel.Calculate()
